# TokenTypes.xlsx edit:
#  - Add a new "context" / CONTEXT keyword row (just above the "data" row).
#  - Add a new "Tag" row for the `#abc_def` tag literal (just above the
#    "TemplateLiteral" row), replacing the concept of the old
#    "UserDefinedKeyWord" (`#abc`) row, which is removed.
#  - Remove the old "UserDefinedKeyWord" (`#abc`) row entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Remove the obsolete "UserDefinedKeyWord" (`#abc`) row first, since
#    it is the highest-numbered row of the three locations being
#    touched; deleting it does not disturb the row numbers of the
#    other two edits (which are both above it).
# ------------------------------------------------------------------
$udkCell = $ws.Columns.Item(2).Find("UserDefinedKeyWord")
if ($udkCell -ne $null) {
    $udkCell.EntireRow.Delete()
}

# ------------------------------------------------------------------
# 2) Insert the new "Tag" row directly above the "TemplateLiteral" row.
# ------------------------------------------------------------------
$tplCell = $ws.Columns.Item(2).Find("TemplateLiteral")
$tplRow = $tplCell.Row
$ws.Rows.Item($tplRow).EntireRow.Insert()

$ws.Range("A" + $tplRow).Value2 = "#abc_def"
$ws.Range("B" + $tplRow).Value2 = "Tag"
$ws.Range("C" + $tplRow).Value2 = "Tag"
$ws.Range("E" + $tplRow).Formula = '=CONCATENATE("  ",B' + $tplRow + ',",    // ",C' + $tplRow + '," `",A' + $tplRow + ',"`    ",D' + $tplRow + ')'

# ------------------------------------------------------------------
# 3) Insert the new "context" row directly above the "data" row.
# ------------------------------------------------------------------
$dataCell = $ws.Columns.Item(1).Find("data")
$dataRow = $dataCell.Row
$ws.Rows.Item($dataRow).EntireRow.Insert()

$ws.Range("A" + $dataRow).Value2 = "context"
$ws.Range("B" + $dataRow).Value2 = "CONTEXT"
$ws.Range("C" + $dataRow).Value2 = "Keyword"
$ws.Range("E" + $dataRow).Formula = '=CONCATENATE("  ",B' + $dataRow + ',",    // ",C' + $dataRow + '," `",A' + $dataRow + ',"`    ",D' + $dataRow + ')'

# ------------------------------------------------------------------
# 4) Restore the view to the top of the sheet with A2 selected (matches
#    the scroll position / selection recorded after the edit).
# ------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("A2").Select()
